$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Preprocessing data" task ---
$ws.Range("A3").Value = "Preprocessing data"
$ws.Range("B3").Value = "19 Oktober 2021 "
$ws.Range("C3").Value = "25 Oktober 2021"
$ws.Range("D3").Value = "25 Oktober 2021"
$ws.Range("E3").Value = "Bibit"
$ws.Range("F3").Value = "Selesai"

# --- Row 4: "Pemilihan Model" task ---
$ws.Range("A4").Value = "Pemilihan Model (Modelling, Model Interface)"
$ws.Range("B4").Value = "31 Oktober 2021"

# Apply the long-date number format before writing the values so the
# custom numFmt is created cleanly against an otherwise-empty cell.
$ws.Range("C4:D4").NumberFormat = "[`$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("C4").Value2 = 44509
$ws.Range("D4").Value2 = 44507

$ws.Range("E4").Value = "Dharma"
$ws.Range("F4").Value = "Selesai"

# --- Column widths: widen A (long task name) and C/D (long date text) ---
$ws.Columns("A").ColumnWidth = 38.053385416666664
$ws.Columns("C:D").ColumnWidth = 15.721354166666666

# --- Selection moves to E5 after data entry ---
$ws.Range("E5").Select()
